$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.335.78'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '2.064.36'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.67'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0760'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '2.370.14'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.776'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.32%  '
$ws.Range("D17").Value = '2.062.82'
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = '37.286.18'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.57%  '
$ws.Range("D21").Value = '0.0₃0809'
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '226.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.42%  '
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.17%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0616'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.55'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.66%  '
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  -3.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '96.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.28%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.463.98'
$ws.Range("E42").Value = '  -0.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0933'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.41%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0212'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("E46").Value = '  -4.47%  '
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '14.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.25%  '
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").Value = '2.256.09'
$ws.Range("E51").Value = '  -0.37%  '
